$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the "Förändrad" (changed) date column (C) for rows 2-20 from 2023-09-06 (45175)
# to 2023-09-14 (45183), keeping existing cell formatting/style intact.
$newDate = Get-Date -Year 2023 -Month 9 -Day 14 -Hour 0 -Minute 0 -Second 0

for ($row = 2; $row -le 20; $row++) {
    $ws.Cells.Item($row, 3).Value = $newDate
}
